$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Cells.Item(2, 5).Value = 61
$ws.Cells.Item(2, 6).Value = 36
$ws.Cells.Item(2, 8).Value = 36
$ws.Cells.Item(3, 6).Value = 20
$ws.Cells.Item(3, 8).Value = 20
$ws.Cells.Item(4, 5).Value = 26
$ws.Cells.Item(4, 6).Value = 13
$ws.Cells.Item(4, 8).Value = 13
$ws.Cells.Item(5, 6).Value = 35
$ws.Cells.Item(5, 8).Value = 35
$ws.Cells.Item(7, 6).Value = 7
$ws.Cells.Item(7, 8).Value = 7
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(9, 6).Value = 3
$ws.Cells.Item(9, 8).Value = 3
$ws.Cells.Item(10, 5).Value = 199
$ws.Cells.Item(10, 6).Value = 85
$ws.Cells.Item(10, 8).Value = 85
$ws.Cells.Item(11, 5).Value = 142
$ws.Cells.Item(11, 6).Value = 74
$ws.Cells.Item(11, 8).Value = 74
$ws.Cells.Item(12, 5).Value = 213
$ws.Cells.Item(12, 6).Value = 108
$ws.Cells.Item(12, 8).Value = 108
$ws.Cells.Item(13, 5).Value = 69
$ws.Cells.Item(13, 6).Value = 33
$ws.Cells.Item(13, 8).Value = 33
$ws.Cells.Item(14, 5).Value = 64
$ws.Cells.Item(14, 6).Value = 27
$ws.Cells.Item(14, 8).Value = 27
$ws.Cells.Item(15, 5).Value = 92
$ws.Cells.Item(15, 6).Value = 24
$ws.Cells.Item(15, 8).Value = 24
$ws.Cells.Item(16, 5).Value = 85
$ws.Cells.Item(16, 6).Value = 41
$ws.Cells.Item(16, 8).Value = 41
$ws.Cells.Item(17, 5).Value = 39
$ws.Cells.Item(20, 5).Value = 52
$ws.Cells.Item(20, 6).Value = 19
$ws.Cells.Item(20, 8).Value = 19
$ws.Cells.Item(21, 5).Value = 66
$ws.Cells.Item(21, 6).Value = 35
$ws.Cells.Item(21, 8).Value = 35
$ws.Cells.Item(22, 5).Value = 84
$ws.Cells.Item(22, 6).Value = 35
$ws.Cells.Item(22, 8).Value = 35
$ws.Cells.Item(23, 5).Value = 98
$ws.Cells.Item(23, 6).Value = 42
$ws.Cells.Item(23, 8).Value = 42
$ws.Cells.Item(24, 5).Value = 96
$ws.Cells.Item(24, 6).Value = 47
$ws.Cells.Item(24, 8).Value = 47
$ws.Cells.Item(25, 5).Value = 90
$ws.Cells.Item(26, 5).Value = 52
$ws.Cells.Item(26, 6).Value = 25
$ws.Cells.Item(26, 8).Value = 25
$ws.Cells.Item(27, 5).Value = 133
$ws.Cells.Item(27, 6).Value = 65
$ws.Cells.Item(27, 8).Value = 65
$ws.Cells.Item(28, 5).Value = 88
$ws.Cells.Item(28, 6).Value = 25
$ws.Cells.Item(28, 8).Value = 25
$ws.Cells.Item(29, 5).Value = 87
$ws.Cells.Item(29, 6).Value = 50
$ws.Cells.Item(29, 8).Value = 50
$ws.Cells.Item(30, 5).Value = 97
$ws.Cells.Item(30, 6).Value = 48
$ws.Cells.Item(30, 8).Value = 48
$ws.Cells.Item(31, 6).Value = 18
$ws.Cells.Item(31, 8).Value = 18
$ws.Cells.Item(32, 5).Value = 93
$ws.Cells.Item(32, 6).Value = 48
$ws.Cells.Item(32, 8).Value = 48
$ws.Cells.Item(33, 5).Value = 121
$ws.Cells.Item(34, 5).Value = 95
$ws.Cells.Item(34, 6).Value = 52
$ws.Cells.Item(34, 8).Value = 52
$ws.Cells.Item(35, 5).Value = 61
$ws.Cells.Item(35, 6).Value = 32
$ws.Cells.Item(35, 8).Value = 32
$ws.Cells.Item(37, 5).Value = 64
$ws.Cells.Item(37, 6).Value = 30
$ws.Cells.Item(37, 8).Value = 30
$ws.Cells.Item(38, 5).Value = 43
$ws.Cells.Item(39, 5).Value = 104
$ws.Cells.Item(39, 6).Value = 37
$ws.Cells.Item(39, 8).Value = 37
$ws.Cells.Item(40, 5).Value = 135
$ws.Cells.Item(40, 6).Value = 54
$ws.Cells.Item(40, 8).Value = 54
$ws.Cells.Item(41, 5).Value = 167
$ws.Cells.Item(41, 6).Value = 61
$ws.Cells.Item(41, 8).Value = 61
$ws.Cells.Item(42, 6).Value = 71
$ws.Cells.Item(42, 8).Value = 71
$ws.Cells.Item(43, 5).Value = 46
$ws.Cells.Item(43, 6).Value = 17
$ws.Cells.Item(43, 8).Value = 17
$ws.Cells.Item(44, 5).Value = 135
$ws.Cells.Item(44, 6).Value = 63
$ws.Cells.Item(44, 8).Value = 63
$ws.Cells.Item(45, 5).Value = 53
$ws.Cells.Item(45, 6).Value = 31
$ws.Cells.Item(45, 8).Value = 31
$ws.Cells.Item(46, 5).Value = 118
$ws.Cells.Item(46, 6).Value = 55
$ws.Cells.Item(46, 8).Value = 55
$ws.Cells.Item(47, 5).Value = 196
$ws.Cells.Item(47, 6).Value = 81
$ws.Cells.Item(47, 8).Value = 81
$ws.Cells.Item(48, 5).Value = 98
$ws.Cells.Item(48, 6).Value = 29
$ws.Cells.Item(48, 8).Value = 29
$ws.Cells.Item(49, 6).Value = 51
$ws.Cells.Item(49, 8).Value = 51
$ws.Cells.Item(50, 5).Value = 95
$ws.Cells.Item(50, 6).Value = 34
$ws.Cells.Item(50, 8).Value = 34
$ws.Cells.Item(51, 5).Value = 95
